$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33, pushing existing rows 33..110 down to 34..111.
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new record.
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = "2023-02-22"
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100103
$ws.Range("H33").Value = "Frutos de hueso (carozo)"
$ws.Range("I33").Value = 100103002
$ws.Range("J33").Value = "Ciruela"
$ws.Range("K33").Value = "Black Amber"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 270
$ws.Range("N33").Value = 10000
$ws.Range("O33").Value = 11000
$ws.Range("P33").Value = 10444
$ws.Range("Q33").Value = "$/bandeja 18 kilos granel"
$ws.Range("R33").Value = "Provincia de Curicó"
$ws.Range("S33").Value = 580
$ws.Range("T33").Value = 18

Write-Host "Done"
